# LOM3005.xlsx update:
#  - "Semestre ideal" value changes from "EM-5" to "EF-5,EM-5" (row 9, cols B & C)
#  - The "Requisitos:" row and its associated value row are removed (rows 22-23)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Semestre ideal" value cells (B9/C9 currently hold "EM-5")
$ws.Cells.Item(9, 2).Value = "EF-5,EM-5"
$ws.Cells.Item(9, 3).Value = "EF-5,EM-5"

# Remove the "Requisitos:" label row and its value row entirely
$ws.Rows("22:23").Delete()
